$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/formatting from H1 into the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 4
